# Fix Training Data Issue (#48)
# The "Date" column (BF) held a malformed label "5-28-2013-14" on every
# data row; due to how NBA stats were shown, it was off by a day and
# should read "2014-05-28" instead. Fix every occurrence in the sheet.
#
# Notes on technique:
# - Values must remain plain text (t="s"/inline string), not get silently
#   auto-converted by Excel into a date serial number, and the fix must
#   not introduce any new cell styles/number formats on the fixed cells
#   (they should keep whatever default style they already had).
# - Directly assigning a date-shaped string via Range.Value (or Range.
#   Value2/Formula as a literal) causes Excel to reinterpret it as a date
#   and stamp a new NumberFormat style on the cell. To avoid that, we
#   stage the literal text as the *result of a formula* in a scratch
#   range (a formula like ="2014-05-28" evaluates to a text value, not a
#   date), copy it, and paste-special just the values into the target
#   cells - this keeps each destination cell's original style untouched.
#   The scratch range is cleared afterwards so it leaves no trace.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "5-28-2013-14"
$newValue = "2014-05-28"
$targetCol = "BF"
$scratchCol = "ZZ"

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

$rowsToFix = @()
for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Range("$targetCol$row")
    if ($cell.Text -eq $oldValue) {
        $rowsToFix += $row
    }
}

foreach ($row in $rowsToFix) {
    $scratch = $ws.Range("$scratchCol$row")
    $scratch.Formula = '="' + $newValue + '"'
    $scratch.Copy()

    $target = $ws.Range("$targetCol$row")
    $target.PasteSpecial(-4163)  # xlPasteValues

    $scratch.ClearContents()
}

$excel.CutCopyMode = 0
